$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Columns E (particip) and F (taxa_sucesso) currently hold fractions (0-1)
# that should instead be expressed as percentages (0-100), i.e. multiplied
# by 100, for rows 2 through 7.
for ($row = 2; $row -le 7; $row++) {
    $eCell = $ws.Cells.Item($row, 5)  # column E
    $fCell = $ws.Cells.Item($row, 6)  # column F

    $eCell.Value2 = $eCell.Value2 * 100
    $fCell.Value2 = $fCell.Value2 * 100
}
